# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the 8 crafting-job sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H51").Value = 7963.1055
$ws.Range("I51").Value = 16771.428
$ws.Range("J51").Value = 2824.9167
$ws.Range("K51").Value = 16771.428
$ws.Range("L51").Value = 2824.9167
$ws.Range("M51").Value = -16287.428
$ws.Range("N51").Value = -3792.9167

$ws.Range("H76").Value = 4366.1665
$ws.Range("I76").Value = 4049.25
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 4049.25
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -3734.25
$ws.Range("N76").Value = -5630

$ws.Range("H79").Value = 4366.1665
$ws.Range("I79").Value = 4049.25
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 4049.25
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -2957.25
$ws.Range("N79").Value = -7184

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 23836.701
$ws.Range("I32").Value = 7155.024
$ws.Range("J32").Value = 163962.8
$ws.Range("K32").Value = 7155.024
$ws.Range("L32").Value = 163962.8
$ws.Range("M32").Value = -6868.024
$ws.Range("N32").Value = -164536.8

$ws.Range("H61").Value = 2132.825
$ws.Range("I61").Value = 1282.1111
$ws.Range("J61").Value = 2828.8635
$ws.Range("K61").Value = 1282.1111
$ws.Range("L61").Value = 2828.8635
$ws.Range("M61").Value = -1070.1111
$ws.Range("N61").Value = -3252.8635

$ws.Range("H110").Value = 71579464
$ws.Range("I110").Value = 71579464
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 71579464
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -71577419
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 10684.55
$ws.Range("I132").Value = 11691.827
$ws.Range("J132").Value = 4137.25
$ws.Range("K132").Value = 35075.481
$ws.Range("L132").Value = 12411.75
$ws.Range("M132").Value = -32545.481
$ws.Range("N132").Value = -17471.75

$ws.Range("H136").Value = 2132.825
$ws.Range("I136").Value = 1282.1111
$ws.Range("J136").Value = 2828.8635
$ws.Range("K136").Value = 3846.3333
$ws.Range("L136").Value = 8486.5905
$ws.Range("M136").Value = -1296.3333
$ws.Range("N136").Value = -13586.5905

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H20").Value = 41281
$ws.Range("I20").Value = 53218.95
$ws.Range("J20").Value = 7172.5713
$ws.Range("K20").Value = 53218.95
$ws.Range("L20").Value = 7172.5713
$ws.Range("M20").Value = -52971.95
$ws.Range("N20").Value = -7666.5713

$ws.Range("H80").Value = 2286.84
$ws.Range("J80").Value = 2896.7646
$ws.Range("L80").Value = 2896.7646
$ws.Range("N80").Value = -4892.7646

$ws.Range("H83").Value = 2286.84
$ws.Range("J83").Value = 2896.7646
$ws.Range("L83").Value = 14483.823
$ws.Range("N83").Value = -24467.823

$ws.Range("H107").Value = 250114460
$ws.Range("I107").Value = 333485700
$ws.Range("J107").Value = 777
$ws.Range("K107").Value = 333485700
$ws.Range("L107").Value = 777
$ws.Range("M107").Value = -333483780
$ws.Range("N107").Value = -4617

$ws.Range("H134").Value = 11182.437
$ws.Range("I134").Value = 12253.8545
$ws.Range("J134").Value = 3835.5715
$ws.Range("K134").Value = 36761.5635
$ws.Range("L134").Value = 11506.7145
$ws.Range("M134").Value = -34226.5635
$ws.Range("N134").Value = -16576.7145

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 35954.215
$ws.Range("I31").Value = 638.94446
$ws.Range("K31").Value = 638.94446
$ws.Range("M31").Value = -343.94446

$ws.Range("H34").Value = 35954.215
$ws.Range("I34").Value = 638.94446
$ws.Range("K34").Value = 638.94446
$ws.Range("M34").Value = -436.94446

$ws.Range("H62").Value = 2466.6667
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3948

$ws.Range("H65").Value = 2466.6667
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -19740

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H117").Value = 7146.579
$ws.Range("I117").Value = 880
$ws.Range("J117").Value = 7494.722
$ws.Range("K117").Value = 2640
$ws.Range("L117").Value = 22484.166
$ws.Range("M117").Value = 802
$ws.Range("N117").Value = -29368.166

$ws.Range("H140").Value = 4941.2905
$ws.Range("I140").Value = 6369
$ws.Range("J140").Value = 2345.4546
$ws.Range("K140").Value = 19107
$ws.Range("L140").Value = 7036.3638
$ws.Range("M140").Value = -13927
$ws.Range("N140").Value = -17396.3638

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 3475.15
$ws.Range("I102").Value = 2191.9
$ws.Range("K102").Value = 2191.9
$ws.Range("M102").Value = -569.9000000000001

$ws.Range("H126").Value = 3372.9092
$ws.Range("I126").Value = 3110.2
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 9330.599999999999
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -6860.599999999999
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 2017.0613
$ws.Range("I132").Value = 1511.9117
$ws.Range("K132").Value = 4535.7351
$ws.Range("M132").Value = -2005.7351

$ws.Range("H136").Value = 20215.23
$ws.Range("J136").Value = 20215.23
$ws.Range("L136").Value = 60645.69
$ws.Range("N136").Value = -65745.69

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 2484.2354
$ws.Range("I7").Value = 1532.1111
$ws.Range("J7").Value = 3555.375
$ws.Range("K7").Value = 1532.1111
$ws.Range("L7").Value = 3555.375
$ws.Range("M7").Value = -1420.1111
$ws.Range("N7").Value = -3779.375

$ws.Range("H40").Value = 46156.695
$ws.Range("I40").Value = 113516
$ws.Range("K40").Value = 113516
$ws.Range("M40").Value = -113380

$ws.Range("H122").Value = 2949.2917
$ws.Range("I122").Value = 2886.75
$ws.Range("J122").Value = 3074.375
$ws.Range("K122").Value = 8660.25
$ws.Range("L122").Value = 9223.125
$ws.Range("M122").Value = -6210.25
$ws.Range("N122").Value = -14123.125

$ws.Range("H126").Value = 2484.2354
$ws.Range("I126").Value = 1532.1111
$ws.Range("J126").Value = 3555.375
$ws.Range("K126").Value = 4596.3333
$ws.Range("L126").Value = 10666.125
$ws.Range("M126").Value = -2126.3333
$ws.Range("N126").Value = -15606.125

$ws.Range("H132").Value = 3745.25
$ws.Range("I132").Value = 3146.8696
$ws.Range("K132").Value = 9440.6088
$ws.Range("M132").Value = -6910.6088

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H122").Value = 1996.4412
$ws.Range("I122").Value = 1523.68
$ws.Range("J122").Value = 3309.6667
$ws.Range("K122").Value = 4571.04
$ws.Range("L122").Value = 9929.000100000001
$ws.Range("M122").Value = -2121.04
$ws.Range("N122").Value = -14829.0001

$ws.Range("H132").Value = 2385.5103
$ws.Range("I132").Value = 2679.4688
$ws.Range("J132").Value = 1832.1765
$ws.Range("K132").Value = 8038.4064
$ws.Range("L132").Value = 5496.529500000001
$ws.Range("M132").Value = -5508.4064
$ws.Range("N132").Value = -10556.5295

